# Generate Report for Handoff
# Adds two newly-tracked localization files (a94391c0-... and d56b4dcd-...)
# as rows 4 and 5 on every sheet: Overview ("File Name"/status table) plus
# the per-locale detail sheets (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$HYPERLINK_UNDERLINE = 2         # xlUnderlineStyleSingle
$HYPERLINK_COLOR     = 15570276  # RGB(0x64,0x95,0xED) == #FF6495ED, matches the workbook's custom HyperLink style
$DATE_FMT            = "yyyy-mm-dd HH:mm:ss"

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = $HYPERLINK_UNDERLINE
    $rng.Font.Color = $HYPERLINK_COLOR
}

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A4").Value2 = "a94391c0-850c-4927-b314-50202ce889ca.md"
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/a94391c0-850c-4927-b314-50202ce889ca.md", "", "", "a94391c0-850c-4927-b314-50202ce889ca.md") | Out-Null
Style-AsHyperlink $ws1.Range("A4")
$ws1.Range("B4").Value2 = "Ready for handoff"
$ws1.Range("C4").Value2 = "Ready for handoff"
$ws1.Range("D4").Value2 = "2016-03-24 02:42:33"
$ws1.Range("D4").NumberFormat = $DATE_FMT

$ws1.Range("A5").Value2 = "d56b4dcd-f316-43e1-b88b-e85df600e8e6.md"
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/d56b4dcd-f316-43e1-b88b-e85df600e8e6.md", "", "", "d56b4dcd-f316-43e1-b88b-e85df600e8e6.md") | Out-Null
Style-AsHyperlink $ws1.Range("A5")
$ws1.Range("B5").Value2 = "Ready for handoff"
$ws1.Range("C5").Value2 = "Ready for handoff"
$ws1.Range("D5").Value2 = "2016-03-24 02:42:33"
$ws1.Range("D5").NumberFormat = $DATE_FMT

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A4").Value2 = "a94391c0-850c-4927-b314-50202ce889ca.md"
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/a94391c0-850c-4927-b314-50202ce889ca.md", "", "", "a94391c0-850c-4927-b314-50202ce889ca.md") | Out-Null
Style-AsHyperlink $ws2.Range("A4")
$ws2.Range("B4").Value2 = ".md"
$ws2.Range("C4").Value2 = "Ready for handoff"
$ws2.Range("D4").Value2 = "a94391c0-850c-4927-b314-50202ce889ca.19275a79a45b9b9e798fe35ffb3a4ec35f5a7bd5.zh-cn.xlf"
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a94391c0-850c-4927-b314-50202ce889ca.19275a79a45b9b9e798fe35ffb3a4ec35f5a7bd5.zh-cn.xlf", "", "", "a94391c0-850c-4927-b314-50202ce889ca.19275a79a45b9b9e798fe35ffb3a4ec35f5a7bd5.zh-cn.xlf") | Out-Null
Style-AsHyperlink $ws2.Range("D4")
$ws2.Range("E4").Value2 = "2016-03-24 02:42:28"
$ws2.Range("E4").NumberFormat = $DATE_FMT
$ws2.Range("H4").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H4").NumberFormat = $DATE_FMT
$ws2.Range("J4").Value2 = "Include"

$ws2.Range("A5").Value2 = "d56b4dcd-f316-43e1-b88b-e85df600e8e6.md"
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/d56b4dcd-f316-43e1-b88b-e85df600e8e6.md", "", "", "d56b4dcd-f316-43e1-b88b-e85df600e8e6.md") | Out-Null
Style-AsHyperlink $ws2.Range("A5")
$ws2.Range("B5").Value2 = ".md"
$ws2.Range("C5").Value2 = "Ready for handoff"
$ws2.Range("D5").Value2 = "d56b4dcd-f316-43e1-b88b-e85df600e8e6.3b859ffac098e06fb64fe26f1e59583dfee5f087.zh-cn.xlf"
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d56b4dcd-f316-43e1-b88b-e85df600e8e6.3b859ffac098e06fb64fe26f1e59583dfee5f087.zh-cn.xlf", "", "", "d56b4dcd-f316-43e1-b88b-e85df600e8e6.3b859ffac098e06fb64fe26f1e59583dfee5f087.zh-cn.xlf") | Out-Null
Style-AsHyperlink $ws2.Range("D5")
$ws2.Range("E5").Value2 = "2016-03-24 02:42:28"
$ws2.Range("E5").NumberFormat = $DATE_FMT
$ws2.Range("H5").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H5").NumberFormat = $DATE_FMT
$ws2.Range("J5").Value2 = "Include"

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A4").Value2 = "a94391c0-850c-4927-b314-50202ce889ca.md"
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/a94391c0-850c-4927-b314-50202ce889ca.md", "", "", "a94391c0-850c-4927-b314-50202ce889ca.md") | Out-Null
Style-AsHyperlink $ws3.Range("A4")
$ws3.Range("B4").Value2 = ".md"
$ws3.Range("C4").Value2 = "Ready for handoff"
$ws3.Range("D4").Value2 = "a94391c0-850c-4927-b314-50202ce889ca.19275a79a45b9b9e798fe35ffb3a4ec35f5a7bd5.de-de.xlf"
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a94391c0-850c-4927-b314-50202ce889ca.19275a79a45b9b9e798fe35ffb3a4ec35f5a7bd5.de-de.xlf", "", "", "a94391c0-850c-4927-b314-50202ce889ca.19275a79a45b9b9e798fe35ffb3a4ec35f5a7bd5.de-de.xlf") | Out-Null
Style-AsHyperlink $ws3.Range("D4")
$ws3.Range("E4").Value2 = "2016-03-24 02:42:33"
$ws3.Range("E4").NumberFormat = $DATE_FMT
$ws3.Range("H4").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H4").NumberFormat = $DATE_FMT
$ws3.Range("J4").Value2 = "Include"

$ws3.Range("A5").Value2 = "d56b4dcd-f316-43e1-b88b-e85df600e8e6.md"
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/d56b4dcd-f316-43e1-b88b-e85df600e8e6.md", "", "", "d56b4dcd-f316-43e1-b88b-e85df600e8e6.md") | Out-Null
Style-AsHyperlink $ws3.Range("A5")
$ws3.Range("B5").Value2 = ".md"
$ws3.Range("C5").Value2 = "Ready for handoff"
$ws3.Range("D5").Value2 = "d56b4dcd-f316-43e1-b88b-e85df600e8e6.3b859ffac098e06fb64fe26f1e59583dfee5f087.de-de.xlf"
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d56b4dcd-f316-43e1-b88b-e85df600e8e6.3b859ffac098e06fb64fe26f1e59583dfee5f087.de-de.xlf", "", "", "d56b4dcd-f316-43e1-b88b-e85df600e8e6.3b859ffac098e06fb64fe26f1e59583dfee5f087.de-de.xlf") | Out-Null
Style-AsHyperlink $ws3.Range("D5")
$ws3.Range("E5").Value2 = "2016-03-24 02:42:33"
$ws3.Range("E5").NumberFormat = $DATE_FMT
$ws3.Range("H5").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H5").NumberFormat = $DATE_FMT
$ws3.Range("J5").Value2 = "Include"
